# Updates a batch of currentAveragePrice / LevePrice / LeveProfit figures
# across several crafting-leve sheets (refreshed market-board snapshot).
$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $val = $Values[$col]
        $cell = $ws.Range("$col$Row")
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

# ALC
Set-Row "ALC" 11 @{ H = 83; I = 83; K = 83; M = 57 }
Set-Row "ALC" 31 @{ H = 3756.889; I = 3756.889; K = 11270.667; M = -11040.667 }
Set-Row "ALC" 33 @{ H = 1849.5; I = 200; K = 200; M = 29 }
Set-Row "ALC" 53 @{ H = 834.1111; I = 801; J = 950; K = 801; L = 950; M = -164; N = -2224 }

# ARM
Set-Row "ARM" 3   @{ H = 400; I = 400; K = 400; M = -285 }
Set-Row "ARM" 13  @{ H = 402; I = 400; K = 400; M = -256 }
Set-Row "ARM" 45  @{ H = 2666.3333; J = 2999.5; L = 2999.5; N = -3753.5 }
Set-Row "ARM" 121 @{ H = 0; J = 0; L = 0; N = $null }
Set-Row "ARM" 122 @{ H = 4426; I = 1250; J = 6014; K = 3750; L = 18042; M = -1300; N = -22942 }

# CRP
Set-Row "CRP" 3 @{ H = 902; I = 902; K = 902; M = -789 }

# CUL
Set-Row "CUL" 34  @{ H = 889.1429000000001; I = 150; J = 1012.3333; K = 450; L = 3036.9999; M = -366; N = -3204.9999 }
Set-Row "CUL" 139 @{ H = 3879.3333; I = 4155.3; J = 2499.5; K = 12465.9; L = 7498.5; M = -7325.900000000001; N = -17778.5 }

# GSM
Set-Row "GSM" 7   @{ H = 803320.6; I = 600; J = 1004000.75; K = 600; L = 1004000.75; M = -488; N = -1004224.75 }
Set-Row "GSM" 8   @{ H = 803320.6; I = 600; J = 1004000.75; K = 600; L = 1004000.75; M = -461; N = -1004278.75 }
Set-Row "GSM" 11  @{ H = 2031814; I = 4417166.5; K = 4417166.5; M = -4417027.5 }
Set-Row "GSM" 20  @{ H = 0; J = 0; L = 0; N = $null }
Set-Row "GSM" 33  @{ H = 10000; J = 10000; L = 10000; N = -10504 }

# LTW
Set-Row "LTW" 14  @{ H = 0; I = 0; K = 0; M = $null }
Set-Row "LTW" 20  @{ H = 15000; J = 15000; L = 15000; N = -15452 }
Set-Row "LTW" 22  @{ H = 2383.3333; I = 1075; K = 1075; M = -780 }
Set-Row "LTW" 27  @{ H = 2383.3333; I = 1075; K = 1075; M = -968 }
Set-Row "LTW" 42  @{ H = 10017997; J = 15023995; L = 15023995; N = -15025121 }
Set-Row "LTW" 47  @{ H = 40000000; J = 0; L = 0; N = $null }
Set-Row "LTW" 49  @{ H = 10017997; J = 15023995; L = 15023995; N = -15024289 }
Set-Row "LTW" 52  @{ H = 40000000; J = 0; L = 0; N = $null }
Set-Row "LTW" 55  @{ H = 643.75; I = 643.75; J = 0; K = 643.75; L = 0; M = -470.75; N = $null }
Set-Row "LTW" 122 @{ H = 4242.8887; I = 4083.8572; K = 12251.5716; M = -9801.571599999999 }
Set-Row "LTW" 137 @{ H = 47200; I = 47200; K = 47200; M = -42100 }

# WVR
Set-Row "WVR" 11  @{ H = 0; I = 0; J = 0; K = 0; L = 0; M = $null; N = $null }
Set-Row "WVR" 19  @{ H = 1000; J = 1000; L = 1000; N = -1348 }
Set-Row "WVR" 38  @{ H = 0; I = 0; K = 0; M = $null }
Set-Row "WVR" 81  @{ H = 450; I = 450; K = 900; M = 161 }
Set-Row "WVR" 84  @{ H = 450; I = 450; K = 4500; M = 804 }
Set-Row "WVR" 126 @{ H = 1000; I = 1000; K = 3000; M = -530 }
Set-Row "WVR" 136 @{ H = 3678.8125; I = 2776; J = 9998.5; K = 8328; L = 29995.5; M = -5778; N = -35095.5 }
